# Assignment14_Workbook.xlsx germanization edit
# - Rename sheet "Agenda Planner" -> "Agenda (Plan)"
# - Translate/adjust several labels and agenda items to German
# - Switch the time columns (A:C, rows 7-19) from 12h AM/PM format to 24h format
# - Update title rich-text run "PLANNER" -> "PLAN" (keeping its formatting)
# - Update the saved selection to D21

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the worksheet (this also updates defined names / formulas that
#     reference the sheet by name automatically) ---
$ws.Name = "Agenda (Plan)"

# --- Header / label cells ---
$ws.Range("D2").Value  = "Team Building Maßnahme"
$ws.Range("C3").Value  = "Ort:"
$ws.Range("C4").Value  = "Datum:"

$ws.Range("B6").Value  = "Ende"
$ws.Range("C6").Value  = "Zeit"
$ws.Range("D6").Value  = "Artikel"
$ws.Range("E6").Value  = "Verantwortlicher"

# --- Agenda item / owner cells ---
$ws.Range("D7").Value  = "Frühstück, Begrüßung"
$ws.Range("D8").Value  = "Einführung"
$ws.Range("D9").Value  = "Übung: Arbeitsbeziehungen"
$ws.Range("D10").Value = "Pause"
$ws.Range("D11").Value = "Wanderung: Cady Falls (Strategie Spiel?)"
$ws.Range("D12").Value = "Mittagessen"
$ws.Range("D13").Value = "Strategie Briefing"
$ws.Range("E14").Value = "Liz nach Info fragen"
$ws.Range("D15").Value = "Übung: Stärken"
$ws.Range("D16").Value = "Pause"
$ws.Range("D17").Value = "Wanderung: Redwoods"
$ws.Range("D18").Value = "Übung: Team Building"
$ws.Range("E18").Value = "Garth, Orgateam"
$ws.Range("D19").Value = "Abendessen"

# --- Footer label ---
$ws.Range("A20").Value = "Gesamt"

# --- Title textbox-like cell (rich text): change only the last run's text
#     ("PLANNER" -> "PLAN") while explicitly re-asserting every run's original
#     font attributes so none of the rich-text formatting is lost on save ---
$titleCell = $ws.Range("E5")

$run1 = $titleCell.Characters(1, 6)     # "AGENDA"
$run1Bold  = $run1.Font.Bold
$run1Color = $run1.Font.Color
$run1Name  = $run1.Font.Name
$run1Size  = $run1.Font.Size

$run2 = $titleCell.Characters(7, 1)     # " "
$run2Bold  = $run2.Font.Bold
$run2Color = $run2.Font.Color
$run2Name  = $run2.Font.Name
$run2Size  = $run2.Font.Size

$run3 = $titleCell.Characters(8, 7)     # "PLANNER"
$run3Bold  = $run3.Font.Bold
$run3Color = $run3.Font.Color
$run3Name  = $run3.Font.Name
$run3Size  = $run3.Font.Size

$run3.Text = "PLAN"

$run1b = $titleCell.Characters(1, 6)
$run1b.Font.Bold  = $run1Bold
$run1b.Font.Color = $run1Color
$run1b.Font.Name  = $run1Name
$run1b.Font.Size  = $run1Size

$run2b = $titleCell.Characters(7, 1)
$run2b.Font.Bold  = $run2Bold
$run2b.Font.Color = $run2Color
$run2b.Font.Name  = $run2Name
$run2b.Font.Size  = $run2Size

$run3b = $titleCell.Characters(8, 4)
$run3b.Font.Bold  = $run3Bold
$run3b.Font.Color = $run3Color
$run3b.Font.Name  = $run3Name
$run3b.Font.Size  = $run3Size

# --- Switch the time grid (A:C, rows 7-19) from 12-hour AM/PM format to a
#     24-hour format, matching the "germanized" locale conventions ---
for ($r = 7; $r -le 19; $r++) {
    $ws.Range("A$r" + ":C$r").NumberFormat = "h:mm;@"
}

# --- Restore/update the saved selection ---
$ws.Range("D21").Select()
